$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert the new "2022-Q3" sheet right before the current "2022-Q2"
#    sheet (currently the 2nd sheet). This shifts every following
#    sheet's position/sheetId down by one, matching the diff.
# ------------------------------------------------------------------
$oldQ2 = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($oldQ2)
$newSheet.Name = "2022-Q3"

# ------------------------------------------------------------------
# 2) Populate the new sheet's header row + single data row, matching
#    the layout used by every other quarterly "fund holdings" sheet.
# ------------------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0

# B column (fund code) keeps leading zeros -> must be stored as text.
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "005395"

$newSheet.Range("C2").Value = "泓德臻远回报灵活配置混合"

# D:G look numeric but are stored as text in every other sheet, so
# force text formatting before assigning them too.
$newSheet.Range("D2:G2").NumberFormat = "@"
$newSheet.Range("D2").Value = "29.76"
$newSheet.Range("E2").Value = "93.35"
$newSheet.Range("F2").Value = "4.34"
$newSheet.Range("G2").Value = "1.2916"

$newSheet.Range("H2").Value = 10

# Reset the number format we forced above back to the sheet default so
# we don't leave a stray "General"-tagged style behind.
$newSheet.Range("B2:G2").Style = "Normal"

# ------------------------------------------------------------------
# 3) Match formatting: header row + the A2 index cell use the same
#    bold/bordered/centered style as every other fund sheet. Grab
#    that style from the sheet that already has it (old "2022-Q2",
#    now shifted to position 3) via copy/paste-special-formats so we
#    reuse the existing style entry instead of inventing a new one.
# ------------------------------------------------------------------
$styleSource = $wb.Worksheets.Item("2022-Q2")
$styleSource.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$styleSource.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 4) Add the corresponding summary row on "总计": a new row 2 for
#    2022-Q3 (pushing the older quarters down, already handled by the
#    later cells shifting automatically via Insert).
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# the freshly inserted row can inherit stray formatting from its
# neighbours; reset it to the sheet default before writing values.
$summary.Range("A2:D2").Style = "Normal"

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 1.29

$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Renumber the running index in column A for every pre-existing row
# (they each move down one row, but keep their own A-value sequence).
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6
$summary.Range("A9").Value = 7

Write-Host "done"
